# Updated cryptos list - refresh Price (col D) and Volume(1h) (col E) figures,
# and shift the coin rows 41-51 up by one (PaxosStandard row dropped, USDD appended)
# to mirror the latest coinranking.com snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table of per-row cell updates: only columns that actually changed are listed.
$rowUpdates = @(
    @{ Row = 2; Cells = @{ "D" = "24.975.85"; "E" = "  -3.88%  " } },
    @{ Row = 3; Cells = @{ "D" = "1.640.94"; "E" = "  -5.74%  " } },
    @{ Row = 4; Cells = @{ "D" = "0.9989" } },
    @{ Row = 5; Cells = @{ "D" = "232.56"; "E" = "  -6.12%  " } },
    @{ Row = 6; Cells = @{ "E" = "  -0.07%  " } },
    @{ Row = 7; Cells = @{ "D" = "0.4750"; "E" = "  -5.42%  " } },
    @{ Row = 8; Cells = @{ "D" = "0.2583"; "E" = "  -5.51%  " } },
    @{ Row = 9; Cells = @{ "D" = "0.06085"; "E" = "  -1.81%  " } },
    @{ Row = 10; Cells = @{ "D" = "0.07027"; "E" = "  -3.24%  " } },
    @{ Row = 11; Cells = @{ "D" = "1.646.79"; "E" = "  -5.21%  " } },
    @{ Row = 12; Cells = @{ "D" = "14.53"; "E" = "  -4.25%  " } },
    @{ Row = 13; Cells = @{ "D" = "0.5849"; "E" = "  -10.59%  " } },
    @{ Row = 14; Cells = @{ "D" = "4.319"; "E" = "  -8.61%  " } },
    @{ Row = 15; Cells = @{ "D" = "73.70"; "E" = "  -5.31%  " } },
    @{ Row = 16; Cells = @{ "E" = "  -0.04%  " } },
    @{ Row = 17; Cells = @{ "D" = "0.9999"; "E" = "  -0.05%  " } },
    @{ Row = 18; Cells = @{ "D" = "24.968.34"; "E" = "  -3.96%  " } },
    @{ Row = 19; Cells = @{ "D" = "0.000006586"; "E" = "  -3.94%  " } },
    @{ Row = 20; Cells = @{ "D" = "11.24"; "E" = "  -5.44%  " } },
    @{ Row = 21; Cells = @{ "D" = "1.857.99"; "E" = "  -5.47%  " } },
    @{ Row = 22; Cells = @{ "D" = "4.313"; "E" = "  -6.53%  " } },
    @{ Row = 23; Cells = @{ "D" = "8.538"; "E" = "  -2.53%  " } },
    @{ Row = 24; Cells = @{ "D" = "5.230"; "E" = "  -3.32%  " } },
    @{ Row = 25; Cells = @{ "D" = "133.47"; "E" = "  -0.70%  " } },
    @{ Row = 26; Cells = @{ "E" = "  -2.53%  " } },
    @{ Row = 27; Cells = @{ "E" = "  -7.91%  " } },
    @{ Row = 28; Cells = @{ "D" = "104.27"; "E" = "  -1.07%  " } },
    @{ Row = 29; Cells = @{ "D" = "1.636"; "E" = "  -8.45%  " } },
    @{ Row = 30; Cells = @{ "D" = "3.899"; "E" = "  -2.50%  " } },
    @{ Row = 31; Cells = @{ "D" = "3.566"; "E" = "  -3.79%  " } },
    @{ Row = 32; Cells = @{ "D" = "0.07562"; "E" = "  -7.43%  " } },
    @{ Row = 33; Cells = @{ "D" = "0.9990"; "E" = "  -0.06%  " } },
    @{ Row = 34; Cells = @{ "D" = "0.04275"; "E" = "  -9.93%  " } },
    @{ Row = 35; Cells = @{ "D" = "2.571"; "E" = "  -3.63%  " } },
    @{ Row = 36; Cells = @{ "D" = "0.5934"; "E" = "  -3.34%  " } },
    @{ Row = 37; Cells = @{ "D" = "0.9286"; "E" = "  -6.88%  " } },
    @{ Row = 38; Cells = @{ "D" = "2.579"; "E" = "  -6.51%  " } },
    @{ Row = 39; Cells = @{ "D" = "0.8919"; "E" = "  +11.31%  " } },
    @{ Row = 40; Cells = @{ "E" = "  -0.09%  " } },
    @{ Row = 41; Cells = @{ "B" = "VeChain"; "C" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; "D" = "0.01496"; "E" = "  -7.67%  " } },
    @{ Row = 42; Cells = @{ "B" = "Quant"; "C" = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; "D" = "98.65"; "E" = "  -2.24%  " } },
    @{ Row = 43; Cells = @{ "B" = "RenderToken"; "C" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; "D" = "1.763"; "E" = "  -9.55%  " } },
    @{ Row = 44; Cells = @{ "B" = "TheSandbox"; "C" = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; "D" = "0.3708"; "E" = "  -5.34%  " } },
    @{ Row = 45; Cells = @{ "B" = "FraxShare"; "C" = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; "D" = "4.662"; "E" = "  -7.01%  " } },
    @{ Row = 46; Cells = @{ "B" = "Algorand"; "C" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; "D" = "0.1102"; "E" = "  -6.30%  " } },
    @{ Row = 47; Cells = @{ "B" = "Aptos"; "C" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; "D" = "6.103"; "E" = "  -4.40%  " } },
    @{ Row = 48; Cells = @{ "B" = "Cronos"; "C" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; "D" = "0.05200"; "E" = "  -1.79%  " } },
    @{ Row = 49; Cells = @{ "B" = "TrueUSD"; "C" = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"; "D" = "1.001"; "E" = "  -0.12%  " } },
    @{ Row = 50; Cells = @{ "B" = "Elrond"; "C" = "https://coinranking.com/coin/omwkOTglq+elrond-egld"; "D" = "28.66"; "E" = "  -7.41%  " } },
    @{ Row = 51; Cells = @{ "B" = "USDD"; "C" = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"; "D" = "0.9991"; "E" = "  +0.15%  " } }
)

foreach ($update in $rowUpdates) {
    $row = $update.Row
    foreach ($col in $update.Cells.Keys) {
        $cellRef = "$col$row"
        $value = $update.Cells[$col]

        $range = $ws.Range($cellRef)
        # Force text storage so values like "24.975.85" or "0.05200" keep their
        # exact original formatting instead of being coerced into numbers,
        # then restore the default (unstyled) cell style so no styling drifts.
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    }
}
